$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1200.1538
$ws.Range("J40").Value = 1289
$ws.Range("L40").Value = 1289
$ws.Range("N40").Value = -1639
$ws.Range("H43").Value = 5129.9
$ws.Range("I43").Value = 2425
$ws.Range("J43").Value = 6933.1665
$ws.Range("K43").Value = 2425
$ws.Range("L43").Value = 6933.1665
$ws.Range("M43").Value = -2356
$ws.Range("N43").Value = -7071.1665
$ws.Range("H62").Value = 3198.1667
$ws.Range("I62").Value = 2297.25
$ws.Range("K62").Value = 2297.25
$ws.Range("M62").Value = -1673.25
$ws.Range("H65").Value = 3198.1667
$ws.Range("I65").Value = 2297.25
$ws.Range("K65").Value = 11486.25
$ws.Range("M65").Value = -8366.25
$ws.Range("H116").Value = 9400.556
$ws.Range("I116").Value = 2533
$ws.Range("K116").Value = 2533
$ws.Range("M116").Value = 909
$ws.Range("H129").Value = 975.75714
$ws.Range("J129").Value = 989.7794
$ws.Range("L129").Value = 2969.3382
$ws.Range("N129").Value = -12969.3382

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1421.8889
$ws.Range("I2").Value = 1049.75
$ws.Range("J2").Value = 1719.6
$ws.Range("K2").Value = 1049.75
$ws.Range("L2").Value = 1719.6
$ws.Range("M2").Value = -936.75
$ws.Range("N2").Value = -1945.6
$ws.Range("H32").Value = 6573.28
$ws.Range("I32").Value = 4197.8687
$ws.Range("J32").Value = 14095.417
$ws.Range("K32").Value = 4197.8687
$ws.Range("L32").Value = 14095.417
$ws.Range("M32").Value = -3910.8687
$ws.Range("N32").Value = -14669.417
$ws.Range("H45").Value = 1483.3334
$ws.Range("I45").Value = 950
$ws.Range("J45").Value = 1750
$ws.Range("K45").Value = 950
$ws.Range("L45").Value = 1750
$ws.Range("M45").Value = -573
$ws.Range("N45").Value = -2504
$ws.Range("H46").Value = 4525
$ws.Range("I46").Value = 1500
$ws.Range("J46").Value = 5533.3335
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 5533.3335
$ws.Range("H74").Value = 2296.9644
$ws.Range("I74").Value = 1837.3
$ws.Range("K74").Value = 1837.3
$ws.Range("M74").Value = -963.3
$ws.Range("H77").Value = 2296.9644
$ws.Range("I77").Value = 1837.3
$ws.Range("K77").Value = 9186.5
$ws.Range("M77").Value = -4818.5
$ws.Range("H116").Value = 1421.8889
$ws.Range("I116").Value = 1049.75
$ws.Range("J116").Value = 1719.6
$ws.Range("K116").Value = 1049.75
$ws.Range("L116").Value = 1719.6
$ws.Range("M116").Value = 1244.25
$ws.Range("N116").Value = -6307.6
$ws.Range("H122").Value = 4405.9473
$ws.Range("I122").Value = 2428.4285
$ws.Range("J122").Value = 5559.5
$ws.Range("K122").Value = 7285.2855
$ws.Range("L122").Value = 16678.5
$ws.Range("M122").Value = -4835.2855
$ws.Range("N122").Value = -21578.5
$ws.Range("H123").Value = 44444
$ws.Range("J123").Value = 44444
$ws.Range("L123").Value = 44444
$ws.Range("N123").Value = -54244
$ws.Range("M46").Value = -1181
$ws.Range("N46").Value = -6171.3335

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1421.8889
$ws.Range("I3").Value = 1049.75
$ws.Range("J3").Value = 1719.6
$ws.Range("K3").Value = 1049.75
$ws.Range("L3").Value = 1719.6
$ws.Range("M3").Value = -935.75
$ws.Range("N3").Value = -1947.6
$ws.Range("H20").Value = 3022
$ws.Range("I20").Value = 3681.4167
$ws.Range("J20").Value = 2230.7
$ws.Range("K20").Value = 3681.4167
$ws.Range("L20").Value = 2230.7
$ws.Range("M20").Value = -3434.4167
$ws.Range("N20").Value = -2724.7
$ws.Range("H25").Value = 2171
$ws.Range("I25").Value = 756.5
$ws.Range("J25").Value = 5000
$ws.Range("K25").Value = 756.5
$ws.Range("L25").Value = 5000
$ws.Range("M25").Value = -521.5
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H94").Value = 2528.8
$ws.Range("I94").Value = 2332
$ws.Range("J94").Value = 3070
$ws.Range("K94").Value = 2332
$ws.Range("L94").Value = 3070
$ws.Range("M94").Value = -1881
$ws.Range("N94").Value = -3972
$ws.Range("N25").Value = -5470

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3851.1052
$ws.Range("I99").Value = 2522.8333
$ws.Range("J99").Value = 6128.143
$ws.Range("K99").Value = 2522.8333
$ws.Range("L99").Value = 6128.143
$ws.Range("M99").Value = -1024.8333
$ws.Range("N99").Value = -9124.143
$ws.Range("H126").Value = 3851.1052
$ws.Range("I126").Value = 2522.8333
$ws.Range("J126").Value = 6128.143
$ws.Range("K126").Value = 7568.499899999999
$ws.Range("L126").Value = 18384.429
$ws.Range("M126").Value = -5098.499899999999
$ws.Range("N126").Value = -23324.429

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2411557.5
$ws.Range("I4").Value = 6027644
$ws.Range("K4").Value = 18082932
$ws.Range("M4").Value = -18082820
$ws.Range("H38").Value = 50.94737
$ws.Range("I38").Value = 27.583334
$ws.Range("K38").Value = 82.75000199999999
$ws.Range("M38").Value = 264.249998
$ws.Range("H68").Value = 7780.533
$ws.Range("I68").Value = 911
$ws.Range("K68").Value = 2733
$ws.Range("M68").Value = -1922
$ws.Range("H71").Value = 7780.533
$ws.Range("I71").Value = 911
$ws.Range("K71").Value = 8199
$ws.Range("M71").Value = -4143
$ws.Range("H113").Value = 647.3889
$ws.Range("I113").Value = 699
$ws.Range("K113").Value = 2097
$ws.Range("M113").Value = 73

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5547.1353
$ws.Range("J70").Value = 5710.8335
$ws.Range("L70").Value = 5710.8335
$ws.Range("N70").Value = -6250.8335
$ws.Range("H73").Value = 5547.1353
$ws.Range("J73").Value = 5710.8335
$ws.Range("L73").Value = 5710.8335
$ws.Range("N73").Value = -7582.8335
$ws.Range("H97").Value = 3126.375
$ws.Range("I97").Value = 2333.3333
$ws.Range("J97").Value = 3602.2
$ws.Range("K97").Value = 2333.3333
$ws.Range("L97").Value = 3602.2
$ws.Range("M97").Value = -1837.3333
$ws.Range("N97").Value = -4594.2
$ws.Range("H126").Value = 4023.077
$ws.Range("I126").Value = 2935.2942
$ws.Range("J126").Value = 5410
$ws.Range("K126").Value = 8805.882599999999
$ws.Range("L126").Value = 16230
$ws.Range("M126").Value = -6335.882599999999
$ws.Range("N126").Value = -21170
$ws.Range("H132").Value = 2926.7407
$ws.Range("I132").Value = 866.6667
$ws.Range("J132").Value = 3515.3333
$ws.Range("K132").Value = 2600.0001
$ws.Range("L132").Value = 10545.9999
$ws.Range("M132").Value = -70.0001000000002
$ws.Range("N132").Value = -15605.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2055.7144
$ws.Range("I93").Value = 965.55554
$ws.Range("J93").Value = 4018
$ws.Range("K93").Value = 965.55554
$ws.Range("L93").Value = 4018
$ws.Range("M93").Value = 282.44446
$ws.Range("N93").Value = -6514
$ws.Range("H100").Value = 2666.6667
$ws.Range("I100").Value = 2400
$ws.Range("J100").Value = 2933.3333
$ws.Range("K100").Value = 2400
$ws.Range("L100").Value = 2933.3333
$ws.Range("M100").Value = -1859
$ws.Range("N100").Value = -4015.3333

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 24726466
$ws.Range("I81").Value = 24726466
$ws.Range("K81").Value = 49452932
$ws.Range("M81").Value = -49451871
$ws.Range("H84").Value = 24726466
$ws.Range("I84").Value = 24726466
$ws.Range("K84").Value = 247264660
$ws.Range("M84").Value = -247259356
$ws.Range("H96").Value = 85427090
$ws.Range("I96").Value = 500250000
$ws.Range("J96").Value = 2462502.5
$ws.Range("K96").Value = 500250000
$ws.Range("L96").Value = 2462502.5
$ws.Range("M96").Value = -500248627
$ws.Range("N96").Value = -2465248.5
$ws.Range("H113").Value = 12893.5
$ws.Range("I113").Value = 100000
$ws.Range("J113").Value = 449.7143
$ws.Range("K113").Value = 300000
$ws.Range("L113").Value = 1349.1429
$ws.Range("M113").Value = -297830
$ws.Range("N113").Value = -5689.1429
$ws.Range("H126").Value = 297577.53
$ws.Range("I126").Value = 1170.2916
$ws.Range("J126").Value = 890392
$ws.Range("K126").Value = 3510.8748
$ws.Range("L126").Value = 2671176
$ws.Range("M126").Value = -1040.8748
$ws.Range("N126").Value = -2676116
